$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "France": append row 17
# ---------------------------------------------------------------------------
$wsFrance = $wb.Worksheets.Item("France")

$wsFrance.Cells.Item(17, 1).Value = 43913.81262007809
$wsFrance.Cells.Item(17, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$wsFrance.Cells.Item(17, 2).Value = 16937
$wsFrance.Cells.Item(17, 3).Value = 676
$wsFrance.Cells.Item(17, 4).Value = 2207
$wsFrance.Cells.Item(17, 2).Resize(1, 3).Style = "Normal"

# ---------------------------------------------------------------------------
# Sheet "Monde": append row 12
# ---------------------------------------------------------------------------
$wsMonde = $wb.Worksheets.Item("Monde")

$wsMonde.Cells.Item(12, 1).Value = 43913.81262007809
$wsMonde.Cells.Item(12, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$wsMonde.Cells.Item(12, 2).Value = 367457
$wsMonde.Cells.Item(12, 3).Value = 16113
$wsMonde.Cells.Item(12, 4).Value = 100879
$wsMonde.Cells.Item(12, 2).Resize(1, 3).Style = "Normal"

# ---------------------------------------------------------------------------
# Sheet "percent": append rows 59-64
# ---------------------------------------------------------------------------
$wsPercent = $wb.Worksheets.Item("percent")

$percentRows = @(
    @(59, "France",          4.61,  4.2,   2.19),
    @(60, "Italie",          17.4,  37.71, 7.37),
    @(61, "Espagne",         9,     13.69, 3.33),
    @(62, "Allemagne",       7.86,  0.73,  0.42),
    @(63, "UK",               1.61,  2.09,  0.14),
    @(64, "Reste du monde",  59.52, 41.58, 86.55)
)

foreach ($rowData in $percentRows) {
    $r = $rowData[0]

    $wsPercent.Cells.Item($r, 1).Value = 43913.81266451103
    $wsPercent.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
    $wsPercent.Cells.Item($r, 2).Value = $rowData[1]
    $wsPercent.Cells.Item($r, 3).Value = $rowData[2]
    $wsPercent.Cells.Item($r, 4).Value = $rowData[3]
    $wsPercent.Cells.Item($r, 5).Value = $rowData[4]
    $wsPercent.Cells.Item($r, 2).Resize(1, 4).Style = "Normal"
}
